$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 87
$ws.Range("E87").Value = 4.378374994833019
$ws.Range("F87").Value = 4.761850487103561
$ws.Range("G87").Value = 0.974451668876546
$ws.Range("H87").Value = 0.1642490337702973
$ws.Range("I87").Value = 3.239674292186175

# Row 88
$ws.Range("E88").Value = -2.797863445655844
$ws.Range("F88").Value = 0.4899402294981954
$ws.Range("G88").Value = -0.6840515148295346
$ws.Range("I88").Value = -2.36794520976924

# Row 89
$ws.Range("E89").Value = -2.912034116707706
$ws.Range("F89").Value = -1.159828578682492
$ws.Range("G89").Value = -0.1304715522525736
$ws.Range("H89").Value = -0.1580530711298803
$ws.Range("I89").Value = -2.623509493325252

# Row 90
$ws.Range("E90").Value = -2.820587636041925
$ws.Range("F90").Value = -1.038027550893114
$ws.Range("G90").Value = -0.1701103953747071
$ws.Range("H90").Value = 0.03569302105791516
$ws.Range("I90").Value = -2.686170261725133

# Row 91
$ws.Range("E91").Value = -6.166515983169825
$ws.Range("F91").Value = -3.674250295393825
$ws.Range("G91").Value = -0.162513297365059
$ws.Range("H91").Value = 0.09133680375806255
$ws.Range("I91").Value = -6.095339489562829

# Row 92
$ws.Range("E92").Value = -2.352309585936044
$ws.Range("F92").Value = -3.562861830463875
$ws.Range("G92").Value = -0.05990943407882802
$ws.Range("H92").Value = -0.01235806542162494
$ws.Range("I92").Value = -2.280042086435591

# Row 93
$ws.Range("E93").Value = -0.6678887022411846
$ws.Range("F93").Value = -3.001825476847245
$ws.Range("G93").Value = -0.03688906968240128
$ws.Range("H93").Value = -0.07352568816414827
$ws.Range("I93").Value = -0.557473944394635
